$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Fall_2015 moved from "future" to "current" status, weight bumped to 0.6
$ws.Range("B2").Value = "current"
$ws.Range("D2").Value = 0.6

# Row 3: Spring_2016 moved from "current" to "past" status, gained a lastmod date,
# weight dropped to 0.3. Copy the date number-format from C4 (an existing "past"
# row) before writing the value so the new cell reuses the same style.
$ws.Range("C4").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C3").Value = (Get-Date -Year 2015 -Month 12 -Day 10 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("B3").Value = "past"
$ws.Range("D3").Value = 0.3

# Weight ripple for the following two "past" rows
$ws.Range("D4").Value = 0.1
$ws.Range("D5").Value = 0.05

# Selection moved to D3 as well
$ws.Range("D3").Select()
